# Fixed errors with MSM implementation.
# The "methodNumberOfLines" sheet incorrectly listed a constructor row for
# each class alongside its "real" method row. Remove the constructor rows
# (SecuritySecureConfig(...) and MallMonitorApplication()), keeping only the
# configure(...) and main(...) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("methodNumberOfLines")

# Remove row 4 first (MallMonitorApplication() / 1) so row indices below it
# aren't affected by the row 2 deletion yet, then remove row 2
# (SecuritySecureConfig(...) / 3).
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(2).Delete()
